# draft-gandhi-spring-stamp-srpm-04.pptx — resize/retitle the
# Session-Sender test-packet slides and tighten the "Figure: ..."
# caption spacing.
#
# EMU -> point conversion is 1 pt = 12700 EMU. The COM host here
# truncates the point->EMU float conversion instead of rounding, so a
# half-EMU epsilon is added to every computed point value to land on
# the exact target EMU value.

$p = $ppt.ActivePresentation
$eps = 0.5/12700

function EMUToPt($emu) {
    return ($emu/12700) + $eps
}

# ---------------------------------------------------------------
# Slide 6 - "Session-Sender Test Packet for Links"
# ---------------------------------------------------------------
$s6 = $p.Slides.Item(6)

# Content Placeholder 2: widen/shift left slightly
$s6cp = $s6.Shapes.Item(3)
$s6cp.Left = EMUToPt 304800
$s6cp.Width = EMUToPt 8572500

# Rectangle 4: Figure caption loses one leading space
$s6r4 = $s6.Shapes.Item(4)
$s6r4tr = $s6r4.TextFrame.TextRange
$s6r4n = $s6r4tr.Paragraphs().Count
$s6r4tr.Paragraphs($s6r4n).Runs(1).Text = "                  Figure: Session-Sender Test Packet"

# ---------------------------------------------------------------
# Slide 7 - "Test Packet for SR-MPLS and SRv6 Policy"
# ---------------------------------------------------------------
$s7 = $p.Slides.Item(7)

# Title 1: reposition, widen, shrink font, rename
$s7Title = $s7.Shapes.Item(1)
$s7Title.Left = EMUToPt 76200
$s7Title.Top = EMUToPt 285750
$s7Title.Width = EMUToPt 4630554
$s7Title.TextFrame.TextRange.Text = "Session-Sender Test Packet for SR-MPLS and SRv6 Policy"
$s7Title.TextFrame.TextRange.Font.Size = 28

# Rectangle 4 (SR-MPLS packet diagram): reposition
$s7R4 = $s7.Shapes.Item(3)
$s7R4.Left = EMUToPt 4876800
$s7R4.Top = EMUToPt 51490

# Rectangle 4 caption text
$s7R4tr = $s7R4.TextFrame.TextRange
$s7R4n = $s7R4tr.Paragraphs().Count
$s7R4tr.Paragraphs($s7R4n).Runs(1).Text = "  Figure: Example session-sender test packet for SR-MPLS Policy"

# Content Placeholder 2 (bullet list): resize
$s7CP = $s7.Shapes.Item(4)
$s7CP.Width = EMUToPt 4478154
$s7CP.Height = EMUToPt 2743201

# Rectangle 8 (SRv6 packet diagram): reposition
$s7R8 = $s7.Shapes.Item(6)
$s7R8.Left = EMUToPt 4876800
$s7R8.Top = EMUToPt 1951063

# Rectangle 8 caption text. The new caption wraps to one fewer line at
# this shape's width, and the shape has <a:spAutoFit/>, so the COM host
# auto-shrinks Height when the run text is written; the diff keeps the
# shape's extent unchanged, so restore Height (and Width, defensively)
# right after.
$s7R8tr = $s7R8.TextFrame.TextRange
$s7R8n = $s7R8tr.Paragraphs().Count
$s7R8w = $s7R8.Width
$s7R8h = $s7R8.Height
$s7R8tr.Paragraphs($s7R8n).Runs(1).Text = "   Figure: Example session-sender test packet for SRv6 Policy"
$s7R8.Width = $s7R8w
$s7R8.Height = $s7R8h

# ---------------------------------------------------------------
# Slide 8 - "Session-Reflector Test Packet"
# ---------------------------------------------------------------
$s8 = $p.Slides.Item(8)

# Rectangle 4: Figure caption loses two leading spaces
$s8r4 = $s8.Shapes.Item(4)
$s8r4tr = $s8r4.TextFrame.TextRange
$s8r4n = $s8r4tr.Paragraphs().Count
$s8r4tr.Paragraphs($s8r4n).Runs(1).Text = "                Figure: Session-Reflector Test Packet"
